$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("R2").Value = [double]"-68.14581589863525"
$ws.Range("S2").Value = [double]"-3.020898060584047e-12"
$ws.Range("T2").Value = [double]"-0.3926034062378131"
$ws.Range("U2").Value = [double]"-5.170781728563633e-13"
$ws.Range("V2").Value = [double]"-0.0001764436175133731"
$ws.Range("W2").Value = [double]"2.961904231906723e-13"
$ws.Range("X2").Value = [double]"-7.408153680025089e-08"
$ws.Range("Y2").Value = [double]"2.720427668334274e-13"
$ws.Range("Z2").Value = [double]"-3.257876712706223e-11"
$ws.Range("AA2").Value = [double]"-1.047500611000503e-12"
$ws.Range("AB2").Value = [double]"8.878655506429485e-14"
$ws.Range("AC2").Value = [double]"8.354821097534107e-13"
$ws.Range("AD2").Value = [double]"-6.702843280855789e-13"
$ws.Range("AE2").Value = [double]"1.265817733494844e-12"
$ws.Range("AF2").Value = [double]"1.445312900633305e-13"

# Row 3
$ws.Range("R3").Value = [double]"3820.53361600509"
$ws.Range("S3").Value = [double]"-2.342903581552818e-12"
$ws.Range("T3").Value = [double]"-13.43056296579693"
$ws.Range("U3").Value = [double]"1.90026194601375e-12"
$ws.Range("V3").Value = [double]"-56.58381757617198"
$ws.Range("W3").Value = [double]"3.434122951293175e-13"
$ws.Range("X3").Value = [double]"-10.71440132254808"
$ws.Range("Y3").Value = [double]"-7.600284009294592e-13"
$ws.Range("Z3").Value = [double]"-1.451650138641211"
$ws.Range("AA3").Value = [double]"7.681002117623181e-13"
$ws.Range("AB3").Value = [double]"-0.2105746429917907"
$ws.Range("AC3").Value = [double]"7.950523979024575e-13"
$ws.Range("AD3").Value = [double]"-0.8379425932359311"
$ws.Range("AE3").Value = [double]"-3.957146984635534e-14"
$ws.Range("AF3").Value = [double]"0.001623900910231786"

# Row 4
$ws.Range("R4").Value = [double]"10041.22891042461"
$ws.Range("S4").Value = [double]"-2.370740118676887e-13"
$ws.Range("T4").Value = [double]"-15.38491883850872"
$ws.Range("U4").Value = [double]"6.826304500105208e-13"
$ws.Range("V4").Value = [double]"0.2337075053410519"
$ws.Range("W4").Value = [double]"1.08336545256439e-13"
$ws.Range("X4").Value = [double]"2.221544324665985"
$ws.Range("Y4").Value = [double]"-6.650674399728485e-13"
$ws.Range("Z4").Value = [double]"3.229782662952247"
$ws.Range("AA4").Value = [double]"-5.344701741907745e-13"
$ws.Range("AB4").Value = [double]"3.866729056542326"
$ws.Range("AC4").Value = [double]"1.348811853387962e-12"
$ws.Range("AD4").Value = [double]"-2.123102370873716"
$ws.Range("AE4").Value = [double]"-9.855560844507763e-14"
$ws.Range("AF4").Value = [double]"0.4053785506600061"

# Row 5
$ws.Range("R5").Value = [double]"3868.987338277233"
$ws.Range("S5").Value = [double]"-5.612616372746389e-13"
$ws.Range("T5").Value = [double]"-151.6735275080437"
$ws.Range("U5").Value = [double]"-1.22400683800898e-12"
$ws.Range("V5").Value = [double]"-154.2990631783649"
$ws.Range("W5").Value = [double]"2.248537369836055e-13"
$ws.Range("X5").Value = [double]"-45.0957393065293"
$ws.Range("Y5").Value = [double]"-9.179861597414763e-13"
$ws.Range("Z5").Value = [double]"-13.75208624883947"
$ws.Range("AA5").Value = [double]"-8.765724578276022e-14"
$ws.Range("AB5").Value = [double]"-4.773636234195753"
$ws.Range("AC5").Value = [double]"4.672726226891419e-13"
$ws.Range("AD5").Value = [double]"-2.544962765292205"
$ws.Range("AE5").Value = [double]"1.312381801705288e-12"
$ws.Range("AF5").Value = [double]"-0.629090928574151"

# Row 6
$ws.Range("R6").Value = [double]"-65.15495286770917"
$ws.Range("S6").Value = [double]"6.995800346651404e-13"
$ws.Range("T6").Value = [double]"-1.86156742206555"
$ws.Range("U6").Value = [double]"-9.69570051032105e-13"
$ws.Range("V6").Value = [double]"-0.01064381918407512"
$ws.Range("W6").Value = [double]"1.349063820378472e-12"
$ws.Range("X6").Value = [double]"-0.0001093483097323885"
$ws.Range("Y6").Value = [double]"-8.286697816841414e-13"
$ws.Range("Z6").Value = [double]"-1.343989007620402e-06"
$ws.Range("AA6").Value = [double]"1.381745549908545e-12"
$ws.Range("AB6").Value = [double]"-1.826358344427081e-08"
$ws.Range("AC6").Value = [double]"4.18300106686472e-13"
$ws.Range("AD6").Value = [double]"-2.639729641896818e-10"
$ws.Range("AE6").Value = [double]"-6.315290763915076e-13"
$ws.Range("AF6").Value = [double]"-6.550721765202965e-12"
